$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# EFO/DO december 2022 release: bump version/date strings in column E (version)
# Note: set the Experimental Factor Ontology (row 4) version first, then the
# Disease Ontology (row 3) release date, so the shared-string table entries
# land in the same order as the authored workbook.
$ws.Range("E4").Value = "v3.49.0"
$ws.Range("E3").Value = "v2022-12-15"

# Update the active selection to match the authored state
$ws.Range("E3").Select()
